# This script updates the Equity Holdings Comparison sheet with refreshed
# data from the quant engine: adds a new "Status" column (D) and replaces
# the old Nov_2025 column with a new Oct_2025 column, recalculating MoM/QoQ,
# plus re-ordering / adding rows for funds that were fully exited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Status"); this shifts the old D:H into E:I
# and Excel automatically carries the header style (bold/border/centered) along.
$ws.Columns("D:D").Insert()

# The old H1 ("MoM") header keeps its style after the shift; I1 ("QoQ") is new
# and starts out unstyled, so copy the header format from H1 onto it.
$ws.Cells.Item(1, 9).Value = "QoQ"
$ws.Range('H1').Copy()
$ws.Range('I1').PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header labels for the new/renamed columns
$ws.Cells.Item(1, 4).Value = 'Status'
$ws.Cells.Item(1, 7).Value = 'Oct_2025'

# Refreshed data: ISIN, Stock Name, Mutual Fund, Status, Jan_2026, Dec_2025,
# Oct_2025, MoM, QoQ for every holding (rows 2-34)

# Row 2: Aurobindo Pharma Limited
$ws.Cells.Item(2, 1).Value = 'INE406A01037'
$ws.Cells.Item(2, 2).Value = 'Aurobindo Pharma Limited'
$ws.Cells.Item(2, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(2, 4).Value = 'Adding Consistently'
$ws.Cells.Item(2, 5).Value = 7.921843
$ws.Cells.Item(2, 6).Value = 6.98705
$ws.Cells.Item(2, 7).Value = 6.397524
$ws.Cells.Item(2, 8).Value = 0.934793
$ws.Cells.Item(2, 9).Value = 1.524319

# Row 3: Samvardhana Motherson International Ltd
$ws.Cells.Item(3, 1).Value = 'INE775A01035'
$ws.Cells.Item(3, 2).Value = 'Samvardhana Motherson International Ltd'
$ws.Cells.Item(3, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(3, 4).Value = 'Adding Consistently'
$ws.Cells.Item(3, 5).Value = 6.56115
$ws.Cells.Item(3, 6).Value = 6.276131
$ws.Cells.Item(3, 7).Value = 4.877879
$ws.Cells.Item(3, 8).Value = 0.2850189999999992
$ws.Cells.Item(3, 9).Value = 1.683271

# Row 4: Lloyds Metals And Energy Limited
$ws.Cells.Item(4, 1).Value = 'INE281B01032'
$ws.Cells.Item(4, 2).Value = 'Lloyds Metals And Energy Limited'
$ws.Cells.Item(4, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(4, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(4, 5).Value = 6.354931
$ws.Cells.Item(4, 6).Value = 6.849578
$ws.Cells.Item(4, 7).Value = 6.433149
$ws.Cells.Item(4, 8).Value = -0.4946470000000005
$ws.Cells.Item(4, 9).Value = -0.07821800000000056

# Row 5: Larsen & Toubro Limited
$ws.Cells.Item(5, 1).Value = 'INE018A01030'
$ws.Cells.Item(5, 2).Value = 'Larsen & Toubro Limited'
$ws.Cells.Item(5, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(5, 4).Value = 'Adding Consistently'
$ws.Cells.Item(5, 5).Value = 5.980584
$ws.Cells.Item(5, 6).Value = 5.592055
$ws.Cells.Item(5, 7).Value = 5.249986
$ws.Cells.Item(5, 8).Value = 0.3885290000000001
$ws.Cells.Item(5, 9).Value = 0.7305980000000005

# Row 6: Reliance Industries Limited
$ws.Cells.Item(6, 1).Value = 'INE002A01018'
$ws.Cells.Item(6, 2).Value = 'Reliance Industries Limited'
$ws.Cells.Item(6, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(6, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(6, 5).Value = 5.877455
$ws.Cells.Item(6, 6).Value = 10.008666
$ws.Cells.Item(6, 7).Value = 10.516761
$ws.Cells.Item(6, 8).Value = -4.131211
$ws.Cells.Item(6, 9).Value = -4.639306

# Row 7: Adani Power Limited
$ws.Cells.Item(7, 1).Value = 'INE814H01029'
$ws.Cells.Item(7, 2).Value = 'Adani Power Limited'
$ws.Cells.Item(7, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(7, 4).Value = 'Adding Consistently'
$ws.Cells.Item(7, 5).Value = 5.445909
$ws.Cells.Item(7, 6).Value = 5.17921
$ws.Cells.Item(7, 7).Value = 5.437755
$ws.Cells.Item(7, 8).Value = 0.266699
$ws.Cells.Item(7, 9).Value = 0.008154000000000217

# Row 8: Jio Financial Services Limited
$ws.Cells.Item(8, 1).Value = 'INE758E01017'
$ws.Cells.Item(8, 2).Value = 'Jio Financial Services Limited'
$ws.Cells.Item(8, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(8, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(8, 5).Value = 3.638142
$ws.Cells.Item(8, 6).Value = 3.796489
$ws.Cells.Item(8, 7).Value = 3.755834
$ws.Cells.Item(8, 8).Value = -0.1583469999999996
$ws.Cells.Item(8, 9).Value = -0.1176919999999999

# Row 9: HDFC Life Insurance Co Ltd
$ws.Cells.Item(9, 1).Value = 'INE795G01014'
$ws.Cells.Item(9, 2).Value = 'HDFC Life Insurance Co Ltd'
$ws.Cells.Item(9, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(9, 4).Value = 'Adding Consistently'
$ws.Cells.Item(9, 5).Value = 2.919895
$ws.Cells.Item(9, 6).Value = 2.696909
$ws.Cells.Item(9, 7).Value = 2.503406
$ws.Cells.Item(9, 8).Value = 0.2229860000000001
$ws.Cells.Item(9, 9).Value = 0.4164889999999999

# Row 10: K.P.R. Mill Limited
$ws.Cells.Item(10, 1).Value = 'INE930H01031'
$ws.Cells.Item(10, 2).Value = 'K.P.R. Mill Limited'
$ws.Cells.Item(10, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(10, 4).Value = 'Adding'
$ws.Cells.Item(10, 5).Value = 2.616977
$ws.Cells.Item(10, 6).Value = 2.553775
$ws.Cells.Item(10, 7).Value = 2.756162
$ws.Cells.Item(10, 8).Value = 0.06320199999999998
$ws.Cells.Item(10, 9).Value = -0.1391849999999999

# Row 11: Britannia Industries Limited
$ws.Cells.Item(11, 1).Value = 'INE216A01030'
$ws.Cells.Item(11, 2).Value = 'Britannia Industries Limited'
$ws.Cells.Item(11, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(11, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(11, 5).Value = 2.478869
$ws.Cells.Item(11, 6).Value = 3.216406
$ws.Cells.Item(11, 7).Value = 2.960406
$ws.Cells.Item(11, 8).Value = -0.7375370000000001
$ws.Cells.Item(11, 9).Value = -0.4815369999999999

# Row 12: Tata Communications Limited
$ws.Cells.Item(12, 1).Value = 'INE151A01013'
$ws.Cells.Item(12, 2).Value = 'Tata Communications Limited'
$ws.Cells.Item(12, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(12, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(12, 5).Value = 2.22615
$ws.Cells.Item(12, 6).Value = 2.329661
$ws.Cells.Item(12, 7).Value = 2.27651
$ws.Cells.Item(12, 8).Value = -0.1035110000000001
$ws.Cells.Item(12, 9).Value = -0.05035999999999996

# Row 13: Adani Wilmar Limited
$ws.Cells.Item(13, 1).Value = 'INE699H01024'
$ws.Cells.Item(13, 2).Value = 'Adani Wilmar Limited'
$ws.Cells.Item(13, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(13, 4).Value = 'Adding'
$ws.Cells.Item(13, 5).Value = 2.016296
$ws.Cells.Item(13, 6).Value = 2.013451
$ws.Cells.Item(13, 7).Value = 2.217973
$ws.Cells.Item(13, 8).Value = 0.002845000000000208
$ws.Cells.Item(13, 9).Value = -0.2016770000000001

# Row 14: Escorts Kubota Limited
$ws.Cells.Item(14, 1).Value = 'INE042A01014'
$ws.Cells.Item(14, 2).Value = 'Escorts Kubota Limited'
$ws.Cells.Item(14, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(14, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(14, 5).Value = 1.933665
$ws.Cells.Item(14, 6).Value = 2.099768
$ws.Cells.Item(14, 7).Value = 2.032858
$ws.Cells.Item(14, 8).Value = -0.1661030000000001
$ws.Cells.Item(14, 9).Value = -0.09919300000000009

# Row 15: Oracle Financial Services Software Ltd
$ws.Cells.Item(15, 1).Value = 'INE881D01027'
$ws.Cells.Item(15, 2).Value = 'Oracle Financial Services Software Ltd'
$ws.Cells.Item(15, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(15, 4).Value = 'Adding Consistently'
$ws.Cells.Item(15, 5).Value = 1.637482
$ws.Cells.Item(15, 6).Value = 1.46054
$ws.Cells.Item(15, 7).Value = 1.538716
$ws.Cells.Item(15, 8).Value = 0.1769420000000002
$ws.Cells.Item(15, 9).Value = 0.09876600000000013

# Row 16: Adani Green Energy Limited
$ws.Cells.Item(16, 1).Value = 'INE364U01010'
$ws.Cells.Item(16, 2).Value = 'Adani Green Energy Limited'
$ws.Cells.Item(16, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(16, 4).Value = 'Reducing'
$ws.Cells.Item(16, 5).Value = 1.60878
$ws.Cells.Item(16, 6).Value = 1.725262
$ws.Cells.Item(16, 7).Value = 0.0
$ws.Cells.Item(16, 8).Value = -0.116482
$ws.Cells.Item(16, 9).Value = 1.60878

# Row 17: 3M India Limited
$ws.Cells.Item(17, 1).Value = 'INE470A01017'
$ws.Cells.Item(17, 2).Value = '3M India Limited'
$ws.Cells.Item(17, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(17, 4).Value = 'Adding Consistently'
$ws.Cells.Item(17, 5).Value = 1.380572
$ws.Cells.Item(17, 6).Value = 1.269292
$ws.Cells.Item(17, 7).Value = 1.016226
$ws.Cells.Item(17, 8).Value = 0.1112799999999998
$ws.Cells.Item(17, 9).Value = 0.3643459999999998

# Row 18: JSW Infrastructure Limited
$ws.Cells.Item(18, 1).Value = 'INE880J01026'
$ws.Cells.Item(18, 2).Value = 'JSW Infrastructure Limited'
$ws.Cells.Item(18, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(18, 4).Value = 'Adding Consistently'
$ws.Cells.Item(18, 5).Value = 1.295972
$ws.Cells.Item(18, 6).Value = 1.285633
$ws.Cells.Item(18, 7).Value = 1.244422
$ws.Cells.Item(18, 8).Value = 0.01033899999999988
$ws.Cells.Item(18, 9).Value = 0.05154999999999998

# Row 19: Premier Energies Limited
$ws.Cells.Item(19, 1).Value = 'INE0BS701011'
$ws.Cells.Item(19, 2).Value = 'Premier Energies Limited'
$ws.Cells.Item(19, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(19, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(19, 5).Value = 1.184429
$ws.Cells.Item(19, 6).Value = 1.24758
$ws.Cells.Item(19, 7).Value = 1.540266
$ws.Cells.Item(19, 8).Value = -0.06315099999999996
$ws.Cells.Item(19, 9).Value = -0.355837

# Row 20: Life Insurance Corporation Of India
$ws.Cells.Item(20, 1).Value = 'INE0J1Y01017'
$ws.Cells.Item(20, 2).Value = 'Life Insurance Corporation Of India'
$ws.Cells.Item(20, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(20, 4).Value = 'Reducing Consistently'
$ws.Cells.Item(20, 5).Value = 1.150905
$ws.Cells.Item(20, 6).Value = 4.069492
$ws.Cells.Item(20, 7).Value = 4.050602
$ws.Cells.Item(20, 8).Value = -2.918587
$ws.Cells.Item(20, 9).Value = -2.899697

# Row 21: Adani Energy Solutions Limited
$ws.Cells.Item(21, 1).Value = 'INE931S01010'
$ws.Cells.Item(21, 2).Value = 'Adani Energy Solutions Limited'
$ws.Cells.Item(21, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(21, 4).Value = 'Reducing'
$ws.Cells.Item(21, 5).Value = 1.107085
$ws.Cells.Item(21, 6).Value = 1.144499
$ws.Cells.Item(21, 7).Value = 1.04491
$ws.Cells.Item(21, 8).Value = -0.03741399999999984
$ws.Cells.Item(21, 9).Value = 0.06217500000000009

# Row 22: Colgate-Palmolive (India) Ltd
$ws.Cells.Item(22, 1).Value = 'INE259A01022'
$ws.Cells.Item(22, 2).Value = 'Colgate-Palmolive (India) Ltd'
$ws.Cells.Item(22, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(22, 4).Value = 'Fresh Entry'
$ws.Cells.Item(22, 5).Value = 0.989349
$ws.Cells.Item(22, 6).Value = 0.0
$ws.Cells.Item(22, 7).Value = 0.0
$ws.Cells.Item(22, 8).Value = 0.989349
$ws.Cells.Item(22, 9).Value = 0.989349

# Row 23: Procter & Gamble Hygiene & Health Care Limited
$ws.Cells.Item(23, 1).Value = 'INE179A01014'
$ws.Cells.Item(23, 2).Value = 'Procter & Gamble Hygiene & Health Care Limited'
$ws.Cells.Item(23, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(23, 4).Value = 'Adding Consistently'
$ws.Cells.Item(23, 5).Value = 0.534404
$ws.Cells.Item(23, 6).Value = 0.530464
$ws.Cells.Item(23, 7).Value = 0.519294
$ws.Cells.Item(23, 8).Value = 0.003939999999999944
$ws.Cells.Item(23, 9).Value = 0.01510999999999996

# Row 24: SBI Cards & Payment Services Ltd
$ws.Cells.Item(24, 1).Value = 'INE018E01016'
$ws.Cells.Item(24, 2).Value = 'SBI Cards & Payment Services Ltd'
$ws.Cells.Item(24, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(24, 4).Value = 'Reducing'
$ws.Cells.Item(24, 5).Value = 0.415008
$ws.Cells.Item(24, 6).Value = 0.42731
$ws.Cells.Item(24, 7).Value = 0.4144
$ws.Cells.Item(24, 8).Value = -0.01230200000000004
$ws.Cells.Item(24, 9).Value = 0.0006079999999999974

# Row 25: Tata Consultancy Services Limited
$ws.Cells.Item(25, 1).Value = 'INE467B01029'
$ws.Cells.Item(25, 2).Value = 'Tata Consultancy Services Limited'
$ws.Cells.Item(25, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(25, 4).Value = 'Complete Exit'
$ws.Cells.Item(25, 5).Value = 0.0
$ws.Cells.Item(25, 6).Value = 0.0
$ws.Cells.Item(25, 7).Value = 2.84746
$ws.Cells.Item(25, 8).Value = 0.0
$ws.Cells.Item(25, 9).Value = -2.84746

# Row 26: UNITED BREWERIES LIMITED
$ws.Cells.Item(26, 1).Value = 'INE686F01025'
$ws.Cells.Item(26, 2).Value = 'UNITED BREWERIES LIMITED'
$ws.Cells.Item(26, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(26, 4).Value = 'Complete Exit'
$ws.Cells.Item(26, 5).Value = 0.0
$ws.Cells.Item(26, 6).Value = 0.0
$ws.Cells.Item(26, 7).Value = 2.382529
$ws.Cells.Item(26, 8).Value = 0.0
$ws.Cells.Item(26, 9).Value = -2.382529

# Row 27: DLF Limited
$ws.Cells.Item(27, 1).Value = 'INE271C01023'
$ws.Cells.Item(27, 2).Value = 'DLF Limited'
$ws.Cells.Item(27, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(27, 4).Value = 'Complete Exit'
$ws.Cells.Item(27, 5).Value = 0.0
$ws.Cells.Item(27, 6).Value = 3.552037
$ws.Cells.Item(27, 7).Value = 3.71664
$ws.Cells.Item(27, 8).Value = -3.552037
$ws.Cells.Item(27, 9).Value = -3.71664

# Row 28: GMR Airports Limited
$ws.Cells.Item(28, 1).Value = 'INE776C01039'
$ws.Cells.Item(28, 2).Value = 'GMR Airports Limited'
$ws.Cells.Item(28, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(28, 4).Value = 'Complete Exit'
$ws.Cells.Item(28, 5).Value = 0.0
$ws.Cells.Item(28, 6).Value = 2.514308
$ws.Cells.Item(28, 7).Value = 2.597342
$ws.Cells.Item(28, 8).Value = -2.514308
$ws.Cells.Item(28, 9).Value = -2.597342

# Row 29: ITC Limited
$ws.Cells.Item(29, 1).Value = 'INE154A01025'
$ws.Cells.Item(29, 2).Value = 'ITC Limited'
$ws.Cells.Item(29, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(29, 4).Value = 'Complete Exit'
$ws.Cells.Item(29, 5).Value = 0.0
$ws.Cells.Item(29, 6).Value = 3.439933
$ws.Cells.Item(29, 7).Value = 3.412505
$ws.Cells.Item(29, 8).Value = -3.439933
$ws.Cells.Item(29, 9).Value = -3.412505

# Row 30: Aditya Birla Lifestyle Brands Limited
$ws.Cells.Item(30, 1).Value = 'INE14LE01019'
$ws.Cells.Item(30, 2).Value = 'Aditya Birla Lifestyle Brands Limited'
$ws.Cells.Item(30, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(30, 4).Value = 'Complete Exit'
$ws.Cells.Item(30, 5).Value = 0.0
$ws.Cells.Item(30, 6).Value = 0.0
$ws.Cells.Item(30, 7).Value = 0.69383
$ws.Cells.Item(30, 8).Value = 0.0
$ws.Cells.Item(30, 9).Value = -0.69383

# Row 31: IRB Infrastructure Developers Limited
$ws.Cells.Item(31, 1).Value = 'INE821I01022'
$ws.Cells.Item(31, 2).Value = 'IRB Infrastructure Developers Limited'
$ws.Cells.Item(31, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(31, 4).Value = 'Complete Exit'
$ws.Cells.Item(31, 5).Value = 0.0
$ws.Cells.Item(31, 6).Value = 0.0
$ws.Cells.Item(31, 7).Value = 4.245558
$ws.Cells.Item(31, 8).Value = 0.0
$ws.Cells.Item(31, 9).Value = -4.245558

# Row 32: Indus Towers Limited
$ws.Cells.Item(32, 1).Value = 'INE121J01017'
$ws.Cells.Item(32, 2).Value = 'Indus Towers Limited'
$ws.Cells.Item(32, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(32, 4).Value = 'Complete Exit'
$ws.Cells.Item(32, 5).Value = 0.0
$ws.Cells.Item(32, 6).Value = 0.0
$ws.Cells.Item(32, 7).Value = 1.840259
$ws.Cells.Item(32, 8).Value = 0.0
$ws.Cells.Item(32, 9).Value = -1.840259

# Row 33: LIC Housing Finance Ltd
$ws.Cells.Item(33, 1).Value = 'INE115A01026'
$ws.Cells.Item(33, 2).Value = 'LIC Housing Finance Ltd'
$ws.Cells.Item(33, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(33, 4).Value = 'Complete Exit'
$ws.Cells.Item(33, 5).Value = 0.0
$ws.Cells.Item(33, 6).Value = 0.0
$ws.Cells.Item(33, 7).Value = 0.150677
$ws.Cells.Item(33, 8).Value = 0.0
$ws.Cells.Item(33, 9).Value = -0.150677

# Row 34: Lupin Limited
$ws.Cells.Item(34, 1).Value = 'INE326A01037'
$ws.Cells.Item(34, 2).Value = 'Lupin Limited'
$ws.Cells.Item(34, 3).Value = 'quant Large and Mid Cap Fund'
$ws.Cells.Item(34, 4).Value = 'Complete Exit'
$ws.Cells.Item(34, 5).Value = 0.0
$ws.Cells.Item(34, 6).Value = 0.0
$ws.Cells.Item(34, 7).Value = 1.097036
$ws.Cells.Item(34, 8).Value = 0.0
$ws.Cells.Item(34, 9).Value = -1.097036
